$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update the user name in cell A2
$ws.Range("A2").Value = "Ayati Arvind"

# Update the selection on the Users sheet to A2
$ws.Activate()
$ws.Range("A2").Select()
